# Auto-generated Excel COM-interop script
# Applies the Chocobo_Profits market-data refresh described in the commit diff:
# updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) for the
# affected Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 115.14286
$ws.Range("I4").Value = 115.14286
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 115.14286
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1.142859999999999
$ws.Range("N4").ClearContents()
# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 480
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -16
$ws.Range("N18").Value = -1168
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -931
# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 11016953
$ws.Range("I70").Value = 23400.5
$ws.Range("J70").Value = 25675024
$ws.Range("K70").Value = 70201.5
$ws.Range("L70").Value = 77025072
$ws.Range("M70").Value = -69931.5
$ws.Range("N70").Value = -77025612
# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 11016953
$ws.Range("I73").Value = 23400.5
$ws.Range("J73").Value = 25675024
$ws.Range("K73").Value = 70201.5
$ws.Range("L73").Value = 77025072
$ws.Range("M73").Value = -69265.5
$ws.Range("N73").Value = -77026944
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2265.2524
$ws.Range("I138").Value = 1163.2593
$ws.Range("J138").Value = 2678.5
$ws.Range("K138").Value = 3489.7779
$ws.Range("L138").Value = 8035.5
$ws.Range("M138").Value = 1650.2221
$ws.Range("N138").Value = -18315.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1096.3
$ws.Range("I2").Value = 1073.3334
$ws.Range("J2").Value = 1130.75
$ws.Range("K2").Value = 1073.3334
$ws.Range("L2").Value = 1130.75
$ws.Range("M2").Value = -960.3334
$ws.Range("N2").Value = -1356.75
# Row 3 (Leve Item ID 2494)
$ws.Range("H3").Value = 8483.333000000001
$ws.Range("I3").Value = 4001.6667
$ws.Range("J3").Value = 12965
$ws.Range("K3").Value = 4001.6667
$ws.Range("L3").Value = 12965
$ws.Range("M3").Value = -3886.6667
$ws.Range("N3").Value = -13195
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1600.1666
$ws.Range("I74").Value = 1141.2972
$ws.Range("J74").Value = 4995.8
$ws.Range("K74").Value = 1141.2972
$ws.Range("L74").Value = 4995.8
$ws.Range("M74").Value = -267.2972
$ws.Range("N74").Value = -6743.8
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1600.1666
$ws.Range("I77").Value = 1141.2972
$ws.Range("J77").Value = 4995.8
$ws.Range("K77").Value = 5706.486
$ws.Range("L77").Value = 24979
$ws.Range("M77").Value = -1338.486
$ws.Range("N77").Value = -33715
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 929.4375
$ws.Range("I110").Value = 906
$ws.Range("J110").Value = 999.75
$ws.Range("K110").Value = 906
$ws.Range("L110").Value = 999.75
$ws.Range("M110").Value = 1139
$ws.Range("N110").Value = -5089.75
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1096.3
$ws.Range("I116").Value = 1073.3334
$ws.Range("J116").Value = 1130.75
$ws.Range("K116").Value = 1073.3334
$ws.Range("L116").Value = 1130.75
$ws.Range("M116").Value = 1220.6666
$ws.Range("N116").Value = -5718.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1096.3
$ws.Range("I3").Value = 1073.3334
$ws.Range("J3").Value = 1130.75
$ws.Range("K3").Value = 1073.3334
$ws.Range("L3").Value = 1130.75
$ws.Range("M3").Value = -959.3334
$ws.Range("N3").Value = -1358.75
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2649.8333
$ws.Range("I105").Value = 2679.8
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2679.8
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -932.8000000000002
$ws.Range("N105").Value = -5994
# Row 122 (Leve Item ID 34096)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
# Row 126 (Leve Item ID 34398)
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3051.2727
$ws.Range("I134").Value = 1561.375
$ws.Range("J134").Value = 7024.3335
$ws.Range("K134").Value = 4684.125
$ws.Range("L134").Value = 21073.0005
$ws.Range("M134").Value = -2149.125
$ws.Range("N134").Value = -26143.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1399.0769
$ws.Range("I16").Value = 1108
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1108
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -821
$ws.Range("N16").Value = -3574
# Row 42 (Leve Item ID 1847)
$ws.Range("H42").Value = 25006.666
$ws.Range("I42").Value = 5020
$ws.Range("J42").Value = 35000
$ws.Range("K42").Value = 5020
$ws.Range("L42").Value = 35000
$ws.Range("M42").Value = -4427
$ws.Range("N42").Value = -36186
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 4214.2
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 6059
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 6059
$ws.Range("M99").Value = -1102
$ws.Range("N99").Value = -9055
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2428
$ws.Range("I105").Value = 2050
$ws.Range("J105").Value = 2995
$ws.Range("K105").Value = 2050
$ws.Range("L105").Value = 2995
$ws.Range("M105").Value = -303
$ws.Range("N105").Value = -6489
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1399.0769
$ws.Range("I113").Value = 1108
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1108
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1062
$ws.Range("N113").Value = -7340
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 4214.2
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 6059
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 18177
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -23117
# Row 130 (Leve Item ID 34689)
$ws.Range("H130").Value = 43780
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43780
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43780
$ws.Range("N130").Value = -53820
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 4758.4375
$ws.Range("I134").Value = 5378.2607
$ws.Range("J134").Value = 3174.4443
$ws.Range("K134").Value = 16134.7821
$ws.Range("L134").Value = 9523.332900000001
$ws.Range("M134").Value = -13599.7821
$ws.Range("N134").Value = -14593.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 587.5
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 766.6667
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 2300.0001
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -2524.0001
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 8227
$ws.Range("I68").Value = 908.8889
$ws.Range("J68").Value = 21399.6
$ws.Range("K68").Value = 2726.6667
$ws.Range("L68").Value = 64198.8
$ws.Range("M68").Value = -1915.6667
$ws.Range("N68").Value = -65820.79999999999
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 8227
$ws.Range("I71").Value = 908.8889
$ws.Range("J71").Value = 21399.6
$ws.Range("K71").Value = 8180.0001
$ws.Range("L71").Value = 192596.4
$ws.Range("M71").Value = -4124.0001
$ws.Range("N71").Value = -200708.4
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 5321.8887
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5321.8887
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15965.6661
$ws.Range("N80").Value = -17837.6661
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 5321.8887
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5321.8887
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 47896.99830000001
$ws.Range("N83").Value = -57256.99830000001
# Row 141 (Leve Item ID 44076)
$ws.Range("H141").Value = 8361.267
$ws.Range("I141").Value = 7977.375
$ws.Range("J141").Value = 8800
$ws.Range("K141").Value = 23932.125
$ws.Range("L141").Value = 26400
$ws.Range("M141").Value = -18752.125
$ws.Range("N141").Value = -36760

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3153.3572
$ws.Range("I122").Value = 2097.25
$ws.Range("J122").Value = 9490
$ws.Range("K122").Value = 6291.75
$ws.Range("L122").Value = 28470
$ws.Range("M122").Value = -3841.75
$ws.Range("N122").Value = -33370
# Row 124 (Leve Item ID 34247)
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 21 (Leve Item ID 2672)
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5498.8613
$ws.Range("I132").Value = 2011.8235
$ws.Range("J132").Value = 8618.842000000001
$ws.Range("K132").Value = 6035.470499999999
$ws.Range("L132").Value = 25856.526
$ws.Range("M132").Value = -3505.470499999999
$ws.Range("N132").Value = -30916.526

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 7579845
$ws.Range("I132").Value = 5046
$ws.Range("J132").Value = 15876054
$ws.Range("K132").Value = 15138
$ws.Range("L132").Value = 47628162
$ws.Range("M132").Value = -12608
$ws.Range("N132").Value = -47633222

